$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Reference cell whose style (font size 20, no wrap) matches the rest of
# column C ("status" column) so the new cells line up with the existing ones.
$refStyleCell = $ws.Range("C4")

# New "???" status value for row 8 (new unique shared string)
$ws.Range("C8").Value = "???"
$ws.Range("C8").Font.Size = $refStyleCell.Font.Size

# New "DONE" status values
$doneCells = @("C9", "C12", "C13", "C14", "C15", "C16")
foreach ($addr in $doneCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "DONE"
    $cell.Font.Size = $refStyleCell.Font.Size
}

# Update the view: scroll so row 11 is at the top and select C17
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("C17").Select()
